$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $first = $parts[0]
            if ($first -eq "dnasr281@gmail.com" -or $first -eq "system") {
                $rest = $parts[1..($parts.Length - 1)]
                $newParts = $rest + $first
                $newVal = [string]::Join(", ", $newParts)
                $cell.Value2 = $newVal
            }
        }
    }
}
